$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("T9 Scores by Spending")
$ws.Range("B2").Value = 83.45539900855027
$ws.Range("C2").Value = 83.93381405396646
$ws.Range("D2").Value = 93.46009572653237
$ws.Range("E2").Value = 96.61087677671375
$ws.Range("F2").Value = 90.36945874402643
$ws.Range("B3").Value = 81.8998257021498
$ws.Range("C3").Value = 83.15528577020937
$ws.Range("D3").Value = 87.1335376073717
$ws.Range("E3").Value = 92.71820457965273
$ws.Range("F3").Value = 81.41859632428398
$ws.Range("B4").Value = 78.50200163320186
$ws.Range("C4").Value = 81.63626134231335
$ws.Range("D4").Value = 73.46258857734237
$ws.Range("E4").Value = 84.3192605609222
$ws.Range("F4").Value = 62.7782334137728
$ws.Range("B5").Value = 76.99720981240274
$ws.Range("C5").Value = 81.0278425571344
$ws.Range("D5").Value = 66.16481311032456
$ws.Range("E5").Value = 81.13395072128019
$ws.Range("F5").Value = 53.5268548869691

$ws = $wb.Worksheets.Item("T9 Scores by Size")
$ws.Range("B2").Value = 83.82159776422071
$ws.Range("C2").Value = 83.92984341754834
$ws.Range("D2").Value = 93.55022469776569
$ws.Range("E2").Value = 96.09943667320715
$ws.Range("F2").Value = 89.88385340844357
$ws.Range("B3").Value = 83.36120143857568
$ws.Range("C3").Value = 83.87386873887871
$ws.Range("D3").Value = 93.58239833305436
$ws.Range("E3").Value = 96.7326541730898
$ws.Range("F3").Value = 90.55799747596197
$ws.Range("B4").Value = 77.746416511437
$ws.Range("C4").Value = 81.34449272598371
$ws.Range("D4").Value = 69.96336073939453
$ws.Range("E4").Value = 82.7666344526415
$ws.Range("F4").Value = 58.28600304906789

$ws = $wb.Worksheets.Item("T9 Scores by Type")
$ws.Range("B2").Value = 83.4654254168185
$ws.Range("C2").Value = 83.9023147557395
$ws.Range("D2").Value = 93.61001987197841
$ws.Range("E2").Value = 96.55022312941765
$ws.Range("F2").Value = 90.39253262434622
$ws.Range("B3").Value = 76.95673306832398
$ws.Range("C3").Value = 80.96663632734915
$ws.Range("D3").Value = 66.54845257144746
$ws.Range("E3").Value = 80.79906211395057
$ws.Range("F3").Value = 53.67220822778149
